$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header in F1, matching the style of the other header cells (B1:E1)
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "time_taken"

# Add time_taken values for each data row (F2:F10), as plain text (inline strings)
$times = @(
    "2021-10-05 13:42:14.837305",
    "2021-10-05 13:42:14.837316",
    "2021-10-05 13:42:14.837320",
    "2021-10-05 13:42:14.837323",
    "2021-10-05 13:42:14.837326",
    "2021-10-05 13:42:14.837329",
    "2021-10-05 13:42:14.837332",
    "2021-10-05 13:42:14.837335",
    "2021-10-05 13:42:14.837338"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
